$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "800/SUP 1"
$ws.Range("C2").Value = "J2545456"
$ws.Range("D2").Value = "JALAL MED"
$ws.Range("G2").Value = 10
$ws.Range("H2").Value = 8000
$ws.Range("J2").Value = 800
$ws.Range("N2").Value = 0
$ws.Range("O2").Value = 7200

# Row 3
$ws.Range("A3").Value = "901/FES "
$ws.Range("B3").Value = "Direction régionale"
$ws.Range("C3").Value = "J207703"
$ws.Range("D3").Value = "ACHENGLI LAILA"
$ws.Range("G3").Value = 10
$ws.Range("H3").Value = 5000
$ws.Range("J3").Value = 500
$ws.Range("O3").Value = 4500

# Row 4
$ws.Range("A4").Value = "901/LF/FES "
$ws.Range("B4").Value = "Logement de fonction"
$ws.Range("C4").Value = "BJ36877"
$ws.Range("D4").Value = "CHARIJI ABDELLAH"
$ws.Range("G4").Value = 10
$ws.Range("H4").Value = 6000
$ws.Range("J4").Value = 600
$ws.Range("O4").Value = 5400

# Row 5 becomes the totals row (previously row 6), old row 5 and 6 data merges/replaced
$ws.Range("A5").Value = " "
$ws.Range("B5").Value = " "
$ws.Range("C5").Value = " "
$ws.Range("D5").Value = " "
$ws.Range("E5").Value = " "
$ws.Range("F5").Value = " "
$ws.Range("G5").Value = " "
$ws.Range("H5").Value = 19000
$ws.Range("J5").Value = 1900
$ws.Range("N5").Value = 0
$ws.Range("O5").Value = 17100

# Remove row 6 entirely (its totals values were moved into row 5 above,
# shrinking the table from A1:O6 to A1:O5)
$ws.Rows("6").Delete()
